$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 88

$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"

$ws.Cells.Item($row, 4).Value = 45121
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = 100112031
$ws.Cells.Item($row, 7).Value = "Poroto verde"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 50
$ws.Cells.Item($row, 11).Value = 18000
$ws.Cells.Item($row, 12).Value = 18000
$ws.Cells.Item($row, 13).Value = 18000
$ws.Cells.Item($row, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item($row, 15).Value = "Perú"
$ws.Cells.Item($row, 16).Value = 720
$ws.Cells.Item($row, 17).Value = 25
$ws.Cells.Item($row, 18).Value = "Hortaliza"
